$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> D value, E value (only set when present)
$updates = @(
    @{ Row = 2;  D = "68.262.98"; E = "  +4.36%  " },
    @{ Row = 3;  D = "3.627.34";  E = "  +4.81%  " },
    @{ Row = 4;  D = $null;       E = "  +0.01%  " },
    @{ Row = 5;  D = "202.25";    E = "  +12.26%  " },
    @{ Row = 6;  D = "577.18";    E = "  +3.63%  " },
    @{ Row = 7;  D = "0.620";     E = "  +3.92%  " },
    @{ Row = 8;  D = $null;       E = "  +0.17%  " },
    @{ Row = 9;  D = $null;       E = "  +6.75%  " },
    @{ Row = 10; D = "61.27";     E = "  +20.45%  " },
    @{ Row = 11; D = $null;       E = "  +6.70%  " },
    @{ Row = 12; D = $null;       E = "  +14.51%  " },
    @{ Row = 13; D = "10.32";     E = "  +9.86%  " },
    @{ Row = 14; D = "4.201.14";  E = "  +4.85%  " },
    @{ Row = 15; D = "3.621.58";  E = "  +5.04%  " },
    @{ Row = 16; D = "19.48";     E = "  +10.78%  " },
    @{ Row = 17; D = $null;       E = "  +1.83%  " },
    @{ Row = 18; D = "68.100.30"; E = "  +4.61%  " },
    @{ Row = 19; D = "12.40";     E = "  +7.14%  " },
    @{ Row = 20; D = $null;       E = "  +5.14%  " },
    @{ Row = 21; D = "409.61";    E = "  +9.24%  " },
    @{ Row = 22; D = "13.14";     E = "  +24.00%  " },
    @{ Row = 23; D = $null;       E = "  +4.76%  " },
    @{ Row = 24; D = "85.64";     E = "  +4.26%  " },
    @{ Row = 25; D = $null;       E = "  +17.17%  " },
    @{ Row = 26; D = "2.93";      E = "  +5.83%  " },
    @{ Row = 27; D = "12.62";     E = "  +7.21%  " },
    @{ Row = 28; D = "6.14";      E = "  +2.49%  " },
    @{ Row = 29; D = $null;       E = "  +9.35%  " },
    @{ Row = 30; D = "7.80";      E = "  +10.59%  " },
    @{ Row = 31; D = "31.78";     E = "  +5.97%  " },
    @{ Row = 32; D = "680.01";    E = "  +12.07%  " },
    @{ Row = 33; D = "12.25";     E = "  +4.28%  " },
    @{ Row = 34; D = $null;       E = "  +5.58%  " },
    @{ Row = 35; D = "63.89";     E = "  +2.35%  " },
    @{ Row = 36; D = "41.93";     E = "  +4.30%  " },
    @{ Row = 37; D = "0.415";     E = "  +5.80%  " },
    @{ Row = 38; D = $null;       E = "  -0.08%  " },
    @{ Row = 39; D = "0.0₃0768";  E = "  +8.88%  " },
    @{ Row = 40; D = $null;       E = "  +18.82%  " },
    @{ Row = 41; D = $null;       E = "  +6.02%  " },
    @{ Row = 42; D = "3.192.85";  E = "  +10.72%  " },
    @{ Row = 43; D = $null;       E = "  -0.09%  " },
    @{ Row = 44; D = $null;       E = "  +12.28%  " },
    @{ Row = 45; D = "2.88";      E = "  +28.35%  " },
    @{ Row = 46; D = "2.86";      E = "  +18.02%  " },
    @{ Row = 47; D = "0.0416";    E = "  +6.59%  " },
    @{ Row = 48; D = $null;       E = "  +5.11%  " },
    @{ Row = 49; D = "8.77";      E = "  +8.99%  " },
    @{ Row = 50; D = "3.09";      E = "  +0.49%  " },
    @{ Row = 51; D = "139.38";    E = "  +1.74%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $cellE = $ws.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
}
